# Automated map update: remove the resolved case ("-627") from row 69
# and shift the remaining cases up by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(69).Delete()
